$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; E = 2; F = 20251103 },
    @{ Row = 3; E = 2; F = 20251103 },
    @{ Row = 4; E = 2; F = 20251103 },
    @{ Row = 5; E = 10; F = 20251115 },
    @{ Row = 6; E = 2; F = 20251103 },
    @{ Row = 7; E = 10; F = 20251115 },
    @{ Row = 8; E = 2; F = 20251103 },
    @{ Row = 9; E = 10; F = 20251115 },
    @{ Row = 10; E = 2; F = 20251110 },
    @{ Row = 11; E = 2; F = 20251103 },
    @{ Row = 12; E = 10; F = 20251115 },
    @{ Row = 13; E = 2; F = 20251103 },
    @{ Row = 14; E = 2; F = 20251103 },
    @{ Row = 15; E = 2; F = 20251103 },
    @{ Row = 16; E = 4; F = 20251109 },
    @{ Row = 17; E = 10; F = 20251115 },
    @{ Row = 18; E = 3; F = 20251108 },
    @{ Row = 19; E = 3; F = 20251108 },
    @{ Row = 20; E = 3; F = 20251108 },
    @{ Row = 21; E = 3; F = 20251108 },
    @{ Row = 22; E = 10; F = 20251115 },
    @{ Row = 23; E = 10; F = 20251115 },
    @{ Row = 24; E = 10; F = 20251115 },
    @{ Row = 25; E = 10; F = 20251115 },
    @{ Row = 26; E = 10; F = 20251115 },
    @{ Row = 27; E = 3; F = 20251111 },
    @{ Row = 28; E = 3; F = 20251108 },
    @{ Row = 29; E = 3; F = 20251108 },
    @{ Row = 30; E = 3; F = 20251108 },
    @{ Row = 31; E = 3; F = 20251108 },
    @{ Row = 32; E = 3; F = 20251108 },
    @{ Row = 33; E = 3; F = 20251108 },
    @{ Row = 34; E = 3; F = 20251108 },
    @{ Row = 35; E = 3; F = 20251108 },
    @{ Row = 36; E = 10; F = 202510929 },
    @{ Row = 37; E = 3; F = 20251108 },
    @{ Row = 38; E = 3; F = 20251108 },
    @{ Row = 39; E = 3; F = 20251108 },
    @{ Row = 40; E = 2; F = 20251110 },
    @{ Row = 41; E = 2; F = 20251110 },
    @{ Row = 42; E = 3; F = 20251108 },
    @{ Row = 43; E = 10; F = 20251115 },
    @{ Row = 44; E = 2; F = 20251110 },
    @{ Row = 45; E = 10; F = 20251115 },
    @{ Row = 46; E = 2; F = 20251110 },
    @{ Row = 47; E = 3; F = 20251108 },
    @{ Row = 48; E = 2; F = 20251110 },
    @{ Row = 49; E = 3; F = 20251111 },
    @{ Row = 50; E = 8; F = 20251113 },
    @{ Row = 51; E = 8; F = 20251113 },
    @{ Row = 52; E = 8; F = 20251113 },
    @{ Row = 53; E = 8; F = 20251113 },
    @{ Row = 54; E = 8; F = 20251113 },
    @{ Row = 55; E = 8; F = 20251113 },
    @{ Row = 56; E = 8; F = 20251113 },
    @{ Row = 57; E = 8; F = 20251113 },
    @{ Row = 58; E = 2; F = 20251107 },
    @{ Row = 59; E = 2; F = 20251107 },
    @{ Row = 60; E = 2; F = 20251107 },
    @{ Row = 61; E = 3; F = 20251111 },
    @{ Row = 62; E = 2; F = 20251107 },
    @{ Row = 63; E = 2; F = 20251107 },
    @{ Row = 64; E = 2; F = 20251107 },
    @{ Row = 65; E = 3; F = 20251108 },
    @{ Row = 66; E = 3; F = 20251108 },
    @{ Row = 67; E = 3; F = 20251108 },
    @{ Row = 68; E = 3; F = 20251108 },
    @{ Row = 69; E = 3; F = 20251108 },
    @{ Row = 70; E = 4; F = 20251109 },
    @{ Row = 71; E = 4; F = 20251109 },
    @{ Row = 72; E = 4; F = 20251109 },
    @{ Row = 73; E = 4; F = 20251109 },
    @{ Row = 74; E = 4; F = 20251109 },
    @{ Row = 75; E = 4; F = 20251109 },
    @{ Row = 76; E = 4; F = 20251109 },
    @{ Row = 77; E = 7; F = 20251112 },
    @{ Row = 78; E = 7; F = 20251112 },
    @{ Row = 79; E = 7; F = 20251112 },
    @{ Row = 80; E = 7; F = 20251112 },
    @{ Row = 81; E = 7; F = 20251112 },
    @{ Row = 82; E = 7; F = 20251112 },
    @{ Row = 83; E = 7; F = 20251112 },
    @{ Row = 84; E = 7; F = 20251112 },
    @{ Row = 85; E = 7; F = 20251112 },
    @{ Row = 86; E = 7; F = 20251112 },
    @{ Row = 87; E = 2; F = 20251110 },
    @{ Row = 88; E = 2; F = 20251110 },
    @{ Row = 89; E = 2; F = 20251110 },
    @{ Row = 90; E = 2; F = 20251110 },
    @{ Row = 91; E = 10; F = 20251115 },
    @{ Row = 92; E = 2; F = 20251110 },
    @{ Row = 93; E = 7; F = 20251112 },
    @{ Row = 94; E = 5; F = 20251113 },
    @{ Row = 95; E = 6; F = 20251111 },
    @{ Row = 96; E = 4; F = 20251109 },
    @{ Row = 97; E = 4; F = 20251109 },
    @{ Row = 98; E = 4; F = 20251109 },
    @{ Row = 99; E = 4; F = 20251109 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    $ws.Cells.Item($u.Row, 6).Value = $u.F
}
